{"js": "// Chapter 4 main text + Table 4.1 update:\n// 1) Turn the \"- description of kegg pathway ... generic)\" paragraph red.\n// 2) Turn the first sentence of the following paragraph red, and split the\n//    \". Also, specify how based on curr\" run after the sentence-ending\n//    \". \" so only the \". \" stays red (rest of the sentence stays black).\n// 3) Mark the two related comments (\"Note that it doesn't have to be too\n//    long!\" and \"Not in detail!\") as resolved/done.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet descriptionPara = null;\nlet fewWordsPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"- description of \") === 0 && t.indexOf(\"pathway\") !== -1) {\n    descriptionPara = paragraphs.items[i];\n  }\n  if (t.indexOf(\"- few words to highlight how some components\") === 0) {\n    fewWordsPara = paragraphs.items[i];\n  }\n}\n\nif (!descriptionPara || !fewWordsPara) {\n  throw new Error(\"Could not locate the target paragraphs\");\n}\n\n// 1) Whole paragraph (all runs + paragraph mark) becomes red.\ndescriptionPara.font.color = \"#FF0000\";\n\n// 2) In the next paragraph, color just the leading sentence red, and color\n// the \". \" that follows it red too (splitting it off from the remainder of\n// the run it currently shares with \"Also, specify how based on curr\").\nconst leadSentence = fewWordsPara.search(\n  \"- few words to highlight how some components are more specific to this pathway and others are more generic (Table 1)\",\n  { matchCase: true }\n);\nconst dotSpace = fewWordsPara.search(\". \", { matchCase: true });\nleadSentence.load(\"text\");\ndotSpace.load(\"text\");\nawait context.sync();\n\nleadSentence.items[0].font.color = \"#FF0000\";\ndotSpace.items[0].font.color = \"#FF0000\";\nawait context.sync();\n\n// 3) Resolve the two comments whose paraId become done=\"1\" in the target:\n// \"Note that it doesn't have to be too long! ...\" and \"Not in detail!\".\nconst comments = context.document.body.getComments();\ncomments.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < comments.items.length; i++) {\n  comments.items[i].load(\"content\");\n}\nawait context.sync();\n\nfor (let i = 0; i < comments.items.length; i++) {\n  const content = comments.items[i].content || \"\";\n  if (\n    content.indexOf(\"Note that it doesn't have to be too long\") !== -1 ||\n    content.indexOf(\"Not in detail\") !== -1\n  ) {\n    comments.items[i].resolved = true;\n  }\n}\nawait context.sync();\n", "ps1": "# Chapter 4 main text + Table 4.1 update:\n# 1) Turn the \"- description of kegg pathway ... generic)\" paragraph red.\n# 2) Turn the first sentence of the following paragraph red, and split the\n#    \". Also, specify how based on curr\" run after the sentence-ending\n#    \". \" so only the \". \" stays red (rest of the sentence stays black).\n# 3) Mark the two related comments (\"Note that it doesn't have to be too\n#    long!\" and \"Not in detail!\") as resolved/done.\n\n$d = $word.ActiveDocument\n\n$wdRed = 255\n\n$descriptionPara = $null\n$fewWordsPara = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $t = $para.Range.Text\n    if ($t.StartsWith(\"- description of \") -and $t -match \"pathway\") {\n        $descriptionPara = $para\n    }\n    if ($t.StartsWith(\"- few words to highlight how some components\")) {\n        $fewWordsPara = $para\n    }\n}\n\n# 1) Whole paragraph (all runs + paragraph mark) becomes red.\n$descriptionPara.Range.Font.Color = $wdRed\n\n# 2) In the next paragraph, color just the leading sentence red.\n$leadRange = $fewWordsPara.Range\n$leadRange.Find.ClearFormatting()\n$leadRange.Find.Text = \"- few words to highlight how some components are more specific to this pathway and others are more generic (Table 1)\"\n$leadRange.Find.Execute() | Out-Null\n$leadRange.Font.Color = $wdRed\n\n# Color the \". \" that follows it red too (splitting it off from the\n# remainder of the run it currently shares with \"Also, specify how based on\n# curr\").\n$dotSpaceRange = $fewWordsPara.Range\n$dotSpaceRange.Find.ClearFormatting()\n$dotSpaceRange.Find.Text = \". \"\n$dotSpaceRange.Find.Execute() | Out-Null\n$dotSpaceRange.Font.Color = $wdRed\n\n# 3) Resolve the two comments whose paraId become done=\"1\" in the target:\n# \"Note that it doesn't have to be too long! ...\" and \"Not in detail!\".\nfor ($i = 1; $i -le $d.Comments.Count; $i++) {\n    $comment = $d.Comments.Item($i)\n    $text = $comment.Range.Text\n    if ($text -like \"Note that it doesn't have to be too long*\" -or $text -like \"Not in detail*\") {\n        $comment.Done = $true\n    }\n}\n"}
